# Applies the "Updated symbol list" edit (Fri Dec 16 14:16:38 UTC 2022) to the
# cryptos worksheet. All cells in columns B:E are stored as text in the
# workbook (prices included), so every numeric-looking value is written as
# text: we force the cell to Text format before assigning it, then restore
# the cell's style so no residual number formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Price-only updates (column D) ---
Set-TextValue "D2"  "248.27"
Set-TextValue "D3"  "24.51"
Set-TextValue "D4"  "5.919"
Set-TextValue "D5"  "0.05896"

# --- Row 6 / Row 7 swap (GateToken <-> KuCoinToken) ---
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D6" "6.578"
$ws.Range("E6").Value = "5KuCoinTokenKCS"

$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D7" "3.423"
$ws.Range("E7").Value = "6GateTokenGT"

# --- More price-only updates ---
Set-TextValue "D8"  "1.334"
Set-TextValue "D9"  "0.7980"
Set-TextValue "D10" "0.1473"
Set-TextValue "D11" "0.07827"
Set-TextValue "D14" "0.09255"
Set-TextValue "D15" "3.560"
Set-TextValue "D16" "0.001664"
Set-TextValue "D17" "0.04752"
Set-TextValue "D18" "0.0006065"
Set-TextValue "D19" "0.006229"
Set-TextValue "D20" "0.005589"
Set-TextValue "D21" "0.001068"
Set-TextValue "D23" "3.700"
Set-TextValue "D24" "2.210"
Set-TextValue "D25" "0.3335"
Set-TextValue "D26" "0.1253"
Set-TextValue "D27" "0.0006479"
Set-TextValue "D40" "0.04391"
Set-TextValue "D41" "0.007016"

# --- Row 42 / Row 43 swap (CEJI <-> BKEXToken) ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1067"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003154"
$ws.Range("E43").Value = "42CEJICEJI"

# --- More price-only updates ---
Set-TextValue "D44" "0.01002"
Set-TextValue "D45" "0.002462"
Set-TextValue "D46" "0.00005890"
Set-TextValue "D48" "0.9908"

# --- Row 49: price + label ---
Set-TextValue "D49" "0.1004"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"

Set-TextValue "D50" "0.00002102"

# --- Row 51: label only ---
$ws.Range("E51").Value = "50SpecialPowerGoldSPG"
